$wb = $excel.ActiveWorkbook

$backlog = $wb.Worksheets.Item("Backlog")
$sprint  = $wb.Worksheets.Item("Sprint-Backlog")

# =====================================================================
# Sheet 1: Backlog
# =====================================================================

$backlog.Range("B1").Value = "Team A"

# Row 6: description (C only) gets the long-form text; name (B) unchanged
$backlog.Range("C6").Value = "Umsetzung des Painters im Framework: Standard TicTacToe. Mit Ausblick auf Animation und Schattierung. Schwarzes Gitter, blaue Kreuze, rote Kreise."

# Row 7 now holds what used to be row 7's rules item, description replaced
$backlog.Range("B7").Value = "Implementierung der TicTacToe-Regeln"
$backlog.Range("C7").Value = "Umsetzung der Rules im Framework: Standard TicTacToe mit 30s Zeitlimit. Mit Ausblick auf Highscore. Gewinn +1 Punkt, Verlieren -2 Punkte."

# Row 8
$backlog.Range("B8").Value = "Implementierung eines menschlichen TicTacToe-Spielers"
$backlog.Range("C8").Value = "Umsetzung eines menschlichen Players im Framework: Standard TicTacToe mit Mausbedingung. Mit Ausblick auf Namen eingeben"

# Row 9
$backlog.Range("B9").Value = "Implementierung eines PC gesteurten TicTacToe-Spielers"
$backlog.Range("C9").Value = "Umsetzung eines PC gesteuerten Players im Framework: Standard TicTacToe."

# Row 10
$backlog.Range("B10").Value = "Erzeugung der TicTacToe-Spieler in einer Fabrik."
$backlog.Range("C10").Value = "Umsetung einer Fabrik zur Erzeugung der menschlichen und PC gesteuerten TicTacToe-Spieler."

# Row 11
$backlog.Range("B11").Value = "Spezifikation eines weiteren Spiels"
$backlog.Range("C11").Value = "Im Entwicklungsteam kann beschlossen werden welches weitere Spiel umgesetzt werden soll, bzw. kann."

# Row 12 used to be empty; now holds the new Highscore backlog item
$backlog.Range("A12").Style = $backlog.Range("A11").Style
$backlog.Range("B12").Style = $backlog.Range("B11").Style
$backlog.Range("C12").Style = $backlog.Range("C11").Style

$backlog.Range("A12").Value = 3
$backlog.Range("B12").Value = "Highscore Umsetzung und Darstellung im Framework Integrieren/Spezifizieren."
$backlog.Range("C12").Value = "Highscore: Wo, wie darstellen? Wo, wie speichern? Wie wo umsetzen?"

# View state: zoom + selection on row 7 (whole row)
$backlog.Activate()
$excel.ActiveWindow.Zoom = 100
$backlog.Range("A7:XFD7").Select()

# =====================================================================
# Sheet 2: Sprint-Backlog
# =====================================================================

# New column D header cell: "Focus-Faktor: 0,5"
$sprint.Range("D1").Style = $sprint.Range("C2").Style
$sprint.Range("D1").Value = "Focus-Faktor: 0,5"

# Column E width (matches sheet1's existing col definition already present)
$sprint.Range("D3").Value = "Due: 04.11.2021"
$sprint.Range("E3").Style = $sprint.Range("D3").Style
$sprint.Range("E3").HorizontalAlignment = -4108
$sprint.Range("E3").VerticalAlignment = -4108
$sprint.Range("E3").WrapText = $false
$sprint.Range("E3").Value = "Done"

$sprint.Range("D4").Value = "Due: 04.11.2021"
$sprint.Range("E4").Style = $sprint.Range("D3").Style
$sprint.Range("E4").HorizontalAlignment = -4108
$sprint.Range("E4").VerticalAlignment = -4108
$sprint.Range("E4").WrapText = $false
$sprint.Range("E4").Value = "Done"

# Row 5: previously fully empty, now filled in
$sprint.Range("A5").Style = $sprint.Range("A4").Style
$sprint.Range("B5").Style = $sprint.Range("B4").Style
$sprint.Range("C5").Style = $sprint.Range("C4").Style
$sprint.Range("D5").Style = $sprint.Range("D4").Style
$sprint.Range("E5").Style = $sprint.Range("D3").Style
$sprint.Range("E5").HorizontalAlignment = -4108
$sprint.Range("E5").VerticalAlignment = -4108
$sprint.Range("E5").WrapText = $false

$sprint.Range("A5").Value = 2
$sprint.Range("B5").Value = "Spezifikation einer Umsetzungsidee für das Spiel TicTacToe"
$sprint.Range("C5").Value = "Spezifikation: Wie soll das TicTacToe aussehen? Wie soll das Aussehen technisch erreicht werden? Wie sollen Regeln umgesetzt werden? Soll es Sonderregeln geben, bzw. Sonderspielfelder, etc.? Welche Spieler (PC, ...) soll es geben? Etc.?"
$sprint.Range("D5").Value = "10min"
$sprint.Range("E5").Value = "Done"

# Row 6: previously fully empty, now filled in
$sprint.Range("A6").Style = $sprint.Range("A4").Style
$sprint.Range("B6").Style = $sprint.Range("B4").Style
$sprint.Range("C6").Style = $sprint.Range("C4").Style
$sprint.Range("D6").Style = $sprint.Range("D4").Style

$sprint.Range("A6").Value = 2
$sprint.Range("B6").Value = "Implementierung der TicTacToe-Darstellung"
$sprint.Range("C6").Value = "Umsetzung des Painters im Framework: Standard TicTacToe. Mit Ausblick auf Animation und Schattierung. Schwarzes Gitter, blaue Kreuze, rote Kreise."
$sprint.Range("D6").Value = "180min"

# Row 7: previously fully empty, now filled in
$sprint.Range("A7").Style = $sprint.Range("A4").Style
$sprint.Range("B7").Style = $sprint.Range("B4").Style
$sprint.Range("C7").Style = $sprint.Range("C4").Style
$sprint.Range("D7").Style = $sprint.Range("D4").Style

$sprint.Range("A7").Value = 2
$sprint.Range("B7").Value = "Implementierung der TicTacToe-Regeln"
$sprint.Range("C7").Value = "Umsetzung der Rules im Framework: Standard TicTacToe mit 30s Zeitlimit. Mit Ausblick auf Highscore. Gewinn +1 Punkt, Verlieren -2 Punkte."
$sprint.Range("D7").Value = "180min"

# View state: tab selected, top-left scrolled to row 2, active cell F6
$sprint.Activate()
$excel.ActiveWindow.ScrollRow = 2
$sprint.Range("F6").Select()
